## Git_Notes.docx - "Added info to clone github repo in new machine"
##
## Appends, after the final paragraph ("Doing pull before push is good
## practice..."), a short new section:
##   - a blank paragraph
##   - a heading-style paragraph "How to Access Github repo between systems"
##   - four new bulleted/numbered ("1)", "2)", ...) ListParagraph items
##     describing how to clone the repo on a second machine.
## The four list items form a brand-new numbered list (its own numId),
## independent of (and restarting from 1 relative to) the existing
## "git clone / git pull" list earlier in the document.

$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------
# 1) Locate the insertion point: end of the document (right after the
#    last paragraph, "Doing pull before push is good practice...").
# ---------------------------------------------------------------------
$anchorText = "Doing pull before push is good practice. First we should get the latest code from remote repo and then push our changes"
$anchor = $d.Content
$found = $anchor.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$docEnd = $d.Content.End
$insertionPoint = $d.Range($docEnd, $docEnd)

# ---------------------------------------------------------------------
# 2) Build the new paragraphs as a single WordprocessingML fragment and
#    drop them in right after the last paragraph. The four list items
#    are plain ListParagraph paragraphs for now; the actual list
#    numbering (numPr/numId) gets attached in step 4 via ListFormat so
#    that Word mints a fresh, correctly-linked numbering definition.
# ---------------------------------------------------------------------
$listParaPPr = "<w:pPr><w:pStyle w:val='ListParagraph'/></w:pPr>"

$fragment = ""
$fragment += "<w:p $wNs/>"
$fragment += "<w:p $wNs><w:r><w:t>How to Access Github repo between systems</w:t></w:r></w:p>"
$fragment += "<w:p $wNs>$listParaPPr<w:r><w:t>Github repo should be public</w:t></w:r></w:p>"
$fragment += "<w:p $wNs>$listParaPPr<w:r><w:t>Generate SSH keys in &#8220;Target machine&#8221; and add these keys into the github account</w:t></w:r></w:p>"
$fragment += "<w:p $wNs>$listParaPPr<w:r><w:t>Create empty local git repo</w:t></w:r></w:p>"
$fragment += "<w:p $wNs>$listParaPPr<w:r><w:t>Use &#8220;git clone ssh-clone-code-from-github-repo&#8221;</w:t></w:r></w:p>"

$insertionPoint.InsertXML($fragment)

# ---------------------------------------------------------------------
# 3) Grab the "1) 2) 3) ..." list format already used earlier in the
#    document (the "git clone / git pull" list) so the new list matches
#    the same numbering style (decimal, "%1)").
# ---------------------------------------------------------------------
$styleSource = $d.Content
$styleSource.Find.Execute("git clone git@github.com:SalmanAli-W/First-Git-demo.git", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$listTemplate = $styleSource.ListFormat.ListTemplate

# ---------------------------------------------------------------------
# 4) Turn the four newly-added paragraphs into one freshly-numbered list
#    (restarts at 1, independent numId) by applying that list template
#    across the whole range in one shot.
# ---------------------------------------------------------------------
$total = $d.Paragraphs.Count
$firstItem = $d.Paragraphs.Item($total - 3)
$lastItem = $d.Paragraphs.Item($total)
$listRange = $d.Range($firstItem.Range.Start, $lastItem.Range.End)
$listRange.ListFormat.ApplyListTemplate($listTemplate, $false)
